# Apply the edits described in the commit: update wording in A8, fill in
# several "I" column time entries on the "Planning effectif" sheet, and
# bump a few existing values in column B. Formulas in N/B27/M28/N28/I27
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planning effectif")

# --- Text tweak -----------------------------------------------------------
$ws.Range("A8").Value() = "CRUD de la table jeuvideo avec ses genres et ses plateformes et ses contenu pegi"

# --- New / updated time values in column I --------------------------------
# I6 is a brand-new entry; copy H6's number format/border so it matches the
# rest of the row before writing the value into it.
$ws.Range("H6").Copy()
$ws.Range("I6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I6").Value() = 0.024305555555555556

$ws.Range("I12").Value() = 0.010416666666666666
$ws.Range("I21").Value() = 0.1423611111111111
$ws.Range("I24").Value() = 0.1076388888888889
$ws.Range("I25").Value() = 0.020833333333333332
$ws.Range("I26").Value() = 0.027777777777777776

# --- Updated totals in column B --------------------------------------------
$ws.Range("B24").Value() = 0.23611111111111113
$ws.Range("B25").Value() = 0.13194444444444445
$ws.Range("B26").Value() = 0.375

# --- Restore the active selection shown in the saved workbook -------------
$ws.Range("I25").Select()
